$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 17
$ws.Range("L3").Value = 27.76

$ws.Range("K4").Value = 52
$ws.Range("L4").Value = 136.46

$ws.Range("K5").Value = 73
$ws.Range("L5").Value = 178.56
